$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new rows at row 4 (shifts old rows 4-29 down to 6-31)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Step 2: Style the new A4:A5 cells to match column A formatting (bold, centered, bordered)
$rngA = $ws.Range("A4:A5")
$rngA.Font.Bold = $true
$rngA.HorizontalAlignment = -4108
$rngA.VerticalAlignment = -4160
$rngA.Borders.LineStyle = 1
$rngA.Borders.Weight = 2

# Step 3: Fill new row 4 (Holden) and row 5 (Rizzie Spiral)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C4").Value2 = 1.013338333496214
$ws.Range("D4").Value2 = 0.9967441570949774
$ws.Range("E4").Value2 = 1.013338333496214
$ws.Range("F4").Value2 = 0.9965236839955979
$ws.Range("G4").Value2 = 1.007285059841564
$ws.Range("H4").Value2 = 0.9913597462499435
$ws.Range("I4").Value2 = 0.9965236839955979
$ws.Range("J4").Value2 = 0.9966798794314301
$ws.Range("K4").Value2 = 0.9965236839955979
$ws.Range("L4").Value2 = 0.9967441570949774
$ws.Range("M4").Value2 = 1.005041245295596
$ws.Range("N4").Value2 = 1.005041245295596
$ws.Range("O4").Value2 = 1.005789183477585
$ws.Range("P4").Value2 = 1.002202058195596
$ws.Range("Q4").Value2 = 1.002202058195596
$ws.Range("R4").Value2 = 1.000782464645597
$ws.Range("S4").Value2 = 1.000782464645597
$ws.Range("T4").Value2 = 1.000321810018288
$ws.Range("C5").Value2 = 0.8983377732019843
$ws.Range("D5").Value2 = 1.02383189597336
$ws.Range("E5").Value2 = 0.8983377732019843
$ws.Range("F5").Value2 = 1.028266140371859
$ws.Range("G5").Value2 = 0.9451119904752066
$ws.Range("H5").Value2 = 1.062707098774203
$ws.Range("I5").Value2 = 1.028266140371859
$ws.Range("J5").Value2 = 1.025124693117545
$ws.Range("K5").Value2 = 1.028266140371859
$ws.Range("L5").Value2 = 1.02383189597336
$ws.Range("M5").Value2 = 0.9610848345876721
$ws.Range("N5").Value2 = 0.9610848345876721
$ws.Range("O5").Value2 = 0.9557605532168503
$ws.Range("P5").Value2 = 0.9834786031824012
$ws.Range("Q5").Value2 = 0.9834786031824013
$ws.Range("R5").Value2 = 0.9946754874797659
$ws.Range("S5").Value2 = 0.9946754874797659
$ws.Range("T5").Value2 = 0.997229931985693

# Step 4: Rename "Thomas Hex" to "Matthies Hex" (now located at B11 after the shift)
$ws.Range("B11").Value = "Matthies Hex"

# Step 5: Append two new rows 30 and 31 with updated A/B index labels
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Michael-SNHex"
$rngA2 = $ws.Range("A30:A31")
$rngA2.Font.Bold = $true
$rngA2.HorizontalAlignment = -4108
$rngA2.VerticalAlignment = -4160
$rngA2.Borders.LineStyle = 1
$rngA2.Borders.Weight = 2